$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Retarget the three data tables from the custom "Table_0" style to the
#    built-in table style {58899947-3AF9-45A7-90A7-1CDE37C0FEAB}.
# ---------------------------------------------------------------------------
$targetStyle = "{58899947-3AF9-45A7-90A7-1CDE37C0FEAB}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyle)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Switch the presentation's design from the "Integral" (Red Violet) theme
#    to the standard "Office Theme" colors.
# ---------------------------------------------------------------------------
function HexToRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (standard Office theme palette)
$officeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
    "4472C4", "70AD47", "0563C1", "954F72"
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme
for ($k = 1; $k -le $themeColors.Count; $k++) {
    $themeColors.Item($k).RGB = HexToRgb($officeColors[$k - 1])
}

# Rename the design/theme and color scheme to match "Office Theme" / "Office".
try { $p.Designs.Item(1).Name = "Office Theme" } catch {}
try { $slide.ColorScheme.Name = "Office" } catch {}
try { $p.SlideMaster.Theme.Name = "Office Theme" } catch {}
